$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.099.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.34%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.299.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.36%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("E5").Value = "  -3.64%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.86%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.626"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.82%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.606"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.05%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0912"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.96%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.07%  "

$ws.Range("E13").Value = "  +0.26%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.972"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.57"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.56%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.646.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.44%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.297.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.75%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.053.34"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.47%  "

$ws.Range("E19").Value = "  -4.56%  "

$ws.Range("E20").Value = "  -1.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "74.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.89%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "257.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.79%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.11%  "

$ws.Range("E26").Value = "  +0.44%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.14%  "

$ws.Range("E28").Value = "  +3.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.43%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.70"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.46%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0894"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.77%  "

$ws.Range("E33").Value = "  -5.82%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.82"
$ws.Range("D34").Style = "Normal"

$ws.Range("E35").Value = "  +10.89%  "

$ws.Range("E36").Value = "  -2.28%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.54"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0353"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.63"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.72%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "71.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.68"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.78%  "

$ws.Range("E43").Value = "  -3.23%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.227"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.46%  "

$ws.Range("E45").Value = "  +0.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.59%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "112.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.46%  "

$ws.Range("E48").Value = "  -1.99%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.60%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.48"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.36%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.569.88"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.92%  "
